$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6283387541770935
$ws.Range("B1").Value = 2.179409027099609
$ws.Range("C1").Value = 8.572210311889648
$ws.Range("D1").Value = 1.998711705207825
$ws.Range("E1").Value = 1.171662211418152
